$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new risk entry (R21) as row 22, mirroring the existing table's
# structure (Fase, Risco, Causa, Evento, Consequência, P, I, P*I).
$ws.Range("A22").Value = "Gestão e Fiscalização do Contrato"
$ws.Range("B22").Value = "R21"
$ws.Range("C22").Value = "Falta de Capacitação sobre o Uso de Inteligência Artificial"
$ws.Range("D22").Value = "Oportunidades de melhoria na rotina de contratações do órgão não são aproveitadas devido à inabilidade de compreender e utilizar soluções baseadas em IA."
$ws.Range("E22").Value = "Permanência de processos ineficientes, desperdício de recursos públicos e perda de competitividade na inovação e automação das atividades de contratação."
$ws.Range("F22").Value = 3
$ws.Range("G22").Value = 4
$ws.Range("H22").Value = 12
